$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell F1 - "time_taken", formatted like the other header cells (bold/centered/bordered)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Timestamp values for F2:F21 - plain text values (must not be reinterpreted as dates/numbers)
$timestamps = @(
    "2021-10-05 13:41:42.490932",
    "2021-10-05 13:41:42.490945",
    "2021-10-05 13:41:42.490949",
    "2021-10-05 13:41:42.490952",
    "2021-10-05 13:41:42.490955",
    "2021-10-05 13:41:42.490958",
    "2021-10-05 13:41:42.490961",
    "2021-10-05 13:41:42.490964",
    "2021-10-05 13:41:42.490967",
    "2021-10-05 13:41:42.490970",
    "2021-10-05 13:41:42.490973",
    "2021-10-05 13:41:42.490976",
    "2021-10-05 13:41:42.490979",
    "2021-10-05 13:41:42.490982",
    "2021-10-05 13:41:42.490985",
    "2021-10-05 13:41:42.490988",
    "2021-10-05 13:41:42.490991",
    "2021-10-05 13:41:42.490994",
    "2021-10-05 13:41:42.490997",
    "2021-10-05 13:41:42.491000"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
